# Added first production level email to btsAccounts
#
# - Switches the active tab to the "Production" sheet.
# - Adds a new account row (ID=0, Email, EmailPassword, MLBPassword) using
#   the same "Hyperlink" cell style already used for the other email
#   accounts on the "Test" sheet.
# - Adds a mailto: hyperlink on the new email cell.
# - Leaves the selection on D2, matching the last cell touched.

$wb = $excel.ActiveWorkbook
$testWs = $wb.Worksheets.Item("Test")
$prodWs = $wb.Worksheets.Item("Production")

$email = "faiyam.daft.54@gmail.com"
$password = "sdFgsdfg892m45"

# New account row values.
$prodWs.Range("A2").Value = 0
$prodWs.Range("B2").Value = $email
$prodWs.Range("C2").Value = $password
$prodWs.Range("D2").Value = $password

# Register the mailto: hyperlink on the email cell.
$prodWs.Hyperlinks.Add($prodWs.Range("B2"), "mailto:" + $email)

# Match the "Hyperlink" cell style already used for B2:B5 on the Test
# sheet (copy formats only, so values are untouched).
$testWs.Range("B2").Copy()
$prodWs.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Production becomes the active/selected sheet, with D2 the active cell.
$prodWs.Activate()
$prodWs.Range("D2").Select()
